$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3, shifting existing rows 3-49 down to 4-50.
# This also carries the formatting (including the date style on column D)
# from the row being pushed down, matching the source workbook's behaviour.
$ws.Rows("3:3").Insert()

# Populate the new row 3 with the new data point.
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = (Get-Date -Year 2023 -Month 4 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = 100112044
$ws.Cells.Item(3, 7).Value = "Perejil"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 300
$ws.Cells.Item(3, 11).Value = 1500
$ws.Cells.Item(3, 12).Value = 1500
$ws.Cells.Item(3, 13).Value = 1500
$ws.Cells.Item(3, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(3, 15).Value = "Región del Maule"
$ws.Cells.Item(3, 16).Value = 1500
$ws.Cells.Item(3, 17).Value = 1
$ws.Cells.Item(3, 18).Value = "Hortaliza"
